$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").ClearContents()
$ws.Range("C7").ClearContents()
